$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.732.31"
$ws.Range("E2").Value = "  +5.28%  "

$ws.Range("D3").Value = "2.269.89"
$ws.Range("E3").Value = "  +3.21%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.33"
$ws.Range("E5").Value = "  +2.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.60"
$ws.Range("E7").Value = "  +6.70%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.432"
$ws.Range("E9").Value = "  +8.03%  "

$ws.Range("E10").Value = "  +16.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.55"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.05"
$ws.Range("E12").Value = "  +17.58%  "

$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "2.605.78"
$ws.Range("E14").Value = "  +3.12%  "

$ws.Range("E15").Value = "  +2.64%  "

$ws.Range("E16").Value = "  +5.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.824"
$ws.Range("E17").Value = "  +4.28%  "

$ws.Range("D18").Value = "2.260.21"
$ws.Range("E18").Value = "  +2.64%  "

$ws.Range("D19").Value = "43.595.43"
$ws.Range("E19").Value = "  +5.29%  "

$ws.Range("D20").Value = "0.0₃0994"
$ws.Range("E20").Value = "  +11.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.76"
$ws.Range("E21").Value = "  +2.98%  "

$ws.Range("E22").Value = "  +1.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.86"
$ws.Range("E23").Value = "  +3.29%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("E25").Value = "  +7.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("E26").Value = "  -1.63%  "

$ws.Range("E27").Value = "  +2.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.44"
$ws.Range("E28").Value = "  +2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.98"
$ws.Range("E29").Value = "  +6.66%  "

$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("E31").Value = "  +2.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.81"
$ws.Range("E32").Value = "  +10.58%  "

$ws.Range("E33").Value = "  +1.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0688"
$ws.Range("E34").Value = "  +6.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.08"
$ws.Range("E35").Value = "  +3.26%  "

$ws.Range("E36").Value = "  +2.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.82"
$ws.Range("E37").Value = "  +6.62%  "

$ws.Range("E38").Value = "  +7.85%  "

$ws.Range("E39").Value = "  -0.68%  "

$ws.Range("E40").Value = "  +5.94%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.43"
$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.35"
$ws.Range("E43").Value = "  +6.76%  "

$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.41"
$ws.Range("E44").Value = "  +22.09%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.46"
$ws.Range("E45").Value = "  +4.17%  "

$ws.Range("E46").Value = "  +1.24%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0960"
$ws.Range("E47").Value = "  +0.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.73"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").Value = "1.478.94"
$ws.Range("E49").Value = "  +1.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.34"
$ws.Range("E50").Value = "  +5.35%  "

$ws.Range("E51").Value = "  +1.61%  "
